# Update cryptocurrency price (D) and 1h volume change (E) columns
# to reflect the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.579.14"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.744.23"
$ws.Range("E3").Value = "  +4.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.36"
$ws.Range("E5").Value = "  +3.46%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  +2.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06261"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.743.58"
$ws.Range("E10").Value = "  +4.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07134"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").Value = "  +6.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6225"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.517"
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.58"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.572.88"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006895"
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.72"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.968.26"
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.844"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.369"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.89"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.815"
$ws.Range("E27").Value = "  +5.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.434"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.04"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.002"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.752"
$ws.Range("E31").Value = "  +3.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07877"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04623"
$ws.Range("E33").Value = "  +7.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.619"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6449"
$ws.Range("E35").Value = "  +6.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9987"
$ws.Range("E36").Value = "  +4.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9427"
$ws.Range("E37").Value = "  +3.54%  "
$ws.Range("E38").Value = "  +15.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.998"
$ws.Range("E39").Value = "  +7.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.430"
$ws.Range("E40").Value = "  -6.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E42").Value = "  +18.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01508"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3921"
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1218"
$ws.Range("E45").Value = "  +8.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.759"
$ws.Range("E46").Value = "  +8.54%  "
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.925"
$ws.Range("E48").Value = "  +6.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.71"
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.269"
$ws.Range("E50").Value = "  +5.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3454"
$ws.Range("E51").Value = "  +3.46%  "
